# Rename the "Nokia Cell Phone Checkout" test's testname from
# "NokiaCellPhoneCheckoutTest" to "nokiaCellPhoneCheckoutTest" (lower-case
# leading 'n') on both the RUNMANAGER and DATA sheets (row 5, column A).
#
# The leading apostrophe forces a text/quote-prefixed entry so the cell
# keeps its existing "quote prefix" number format/style (matching the
# original formatting of that cell) instead of Excel resetting it when a
# plain string is assigned.

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsRunManager.Range("A5").Value = "'nokiaCellPhoneCheckoutTest"

$wsData = $wb.Worksheets.Item("DATA")
$wsData.Range("A5").Value = "'nokiaCellPhoneCheckoutTest"
